$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 values
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 2

# Delete row 4 entirely (shrinks used range / dimension to A1:B3)
$ws.Range("A4:B4").Delete()
